$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1075
$ws1.Range("F4").Value = 2481

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 16

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 16
$ws4.Range("F5").Value = 1075
$ws4.Range("F6").Value = 2481
